# Apply cell-value updates to match the target revision of cryptos.xlsx.
# Columns D (Price) and E (Volume(1h)) are stored as text that looks numeric/
# percentage-like, so a leading apostrophe (Excel's text quote-prefix) is used
# to force the value to remain plain text instead of being parsed as a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.54"
$ws.Range("E2").Value = "'0.89%"
$ws.Range("D3").Value = "'29.87"
$ws.Range("E3").Value = "'9.93%"
$ws.Range("D4").Value = "'5.183"
$ws.Range("E4").Value = "'1.28%"
$ws.Range("D5").Value = "'0.05697"
$ws.Range("E5").Value = "'0.55%"
$ws.Range("D6").Value = "'6.593"
$ws.Range("E6").Value = "'1.94%"
$ws.Range("D7").Value = "'0.8585"
$ws.Range("E7").Value = "'4.52%"
$ws.Range("D8").Value = "'0.8750"
$ws.Range("E8").Value = "'3.58%"
$ws.Range("D9").Value = "'0.1365"
$ws.Range("E9").Value = "'3.01%"
$ws.Range("D10").Value = "'0.07086"
$ws.Range("E10").Value = "'2.17%"
$ws.Range("D11").Value = "'0.02868"
$ws.Range("E11").Value = "'-0.58%"
$ws.Range("D12").Value = "'0.09385"
$ws.Range("E12").Value = "'-0.01%"
$ws.Range("D13").Value = "'0.001523"
$ws.Range("E13").Value = "'0.16%"
$ws.Range("D14").Value = "'0.04155"
$ws.Range("E14").Value = "'1.01%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006001"
$ws.Range("E15").Value = "'-0.29%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006027"
$ws.Range("E16").Value = "'-3.09%"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "'0.007491"
$ws.Range("E17").Value = "'5,108.69%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.480"
$ws.Range("E18").Value = "'-0.92%"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'3.047"
$ws.Range("E19").Value = "'1.58%"
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").Value = "'2.186"
$ws.Range("E20").Value = "'-1.84%"
$ws.Range("E21").Value = "'1.00%"
$ws.Range("D22").Value = "'0.03263"
$ws.Range("E22").Value = "'3.03%"
$ws.Range("D23").Value = "'0.1301"
$ws.Range("E23").Value = "'3.64%"
$ws.Range("D24").Value = "'3.478"
$ws.Range("E24").Value = "'-2.28%"
$ws.Range("E25").Value = "'0.45%"
$ws.Range("D26").Value = "'0.005087"
$ws.Range("E26").Value = "'14.23%"
$ws.Range("D27").Value = "'0.001220"
$ws.Range("E27").Value = "'0.01%"
$ws.Range("E28").Value = "'23.50%"
$ws.Range("D40").Value = "'0.03750"
$ws.Range("E40").Value = "'2.08%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1071"
$ws.Range("E41").Value = "'1.71%"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002540"
$ws.Range("E42").Value = "'10.94%"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003457"
$ws.Range("E43").Value = "'-42.94%"
$ws.Range("D44").Value = "'0.009393"
$ws.Range("E44").Value = "'-3.00%"
$ws.Range("D45").Value = "'0.00005114"
$ws.Range("E45").Value = "'-3.82%"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("D47").Value = "'0.07102"
$ws.Range("E47").Value = "'-29.67%"
$ws.Range("D48").Value = "'0.002709"
$ws.Range("E48").Value = "'5.14%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.05%"
